$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1 (Preparation paragraph): re-word the list of items to print, adding
# file extensions / item counts, and splitting "AmdahlLawHandout.docx" out
# with ".pdf" naming.
# ---------------------------------------------------------------------------

$ok1 = $d.Content.Find.Execute(
    "Print out the items given in the pack, one of ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Print out the items given in the pack, ",
    2)
Write-Output "replace1: $ok1"

$ok2 = $d.Content.Find.Execute(
    "AmdahlLawProgram.docx and 1 per student of  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1 of AmdahlLawProgram.pdf, 1 of AmdahlLawBoard ",
    2)
Write-Output "replace2: $ok2"

$ok3 = $d.Content.Find.Execute(
    "^tAmdahlLawHandout.docx",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "^tand 1 per student of AmdahlLawHandout.pdf.",
    2)
Write-Output "replace3: $ok3"

# ---------------------------------------------------------------------------
# Edit 2 (exercises bullet): merge the "AmdahlLaw" / "Handout.docx..." runs
# back into continuous text referencing AmdahlLawHandout.docx (the wording
# itself is unchanged, only the run layout is consolidated).
# ---------------------------------------------------------------------------

$ok4 = $d.Content.Find.Execute(
    "Give students the exercises from AmdahlLawHandout.docx to further explain Amdahl",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Give students the exercises from AmdahlLawHandout.docx to further explain Amdahl",
    2)
Write-Output "replace4: $ok4"

$d.Save()
